# fix academic offer update
# - Correct the name of "Universidad Politécnica de Amozoc" -> "...de Mota"
# - Add the missing "Instituto Tecnológico Superior de la Sierra Norte de
#   Puebla" (Zacatlán) record as a new row at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the existing university name in column A of row 17.
$ws.Range("A17").Value = "Universidad Politécnica de Amozoc de Mota"

# Append the new row (29) with the missing university's data.
$ws.Range("A29").Value = "Instituto Tecnológico Superior de la Sierra Norte de Puebla"
$ws.Range("B29").Value = "http://www.itssnp.edu.mx/v2/"
$ws.Range("C29").Value = "Zacatlán"
$ws.Range("D29").Value = "Puebla"
$ws.Range("E29").Value = "797 975 16 94"
$ws.Range("F29").Value = "Av. José Luis Martínez Vázquez No. 2000. Jicolapa, Zacatlán, Pue"

# Match the author's final cursor position.
$ws.Range("C34").Select()
